$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row 33 (shift down), inheriting cell formatting/styles from
# the row above (row 32), so the new row matches the styling used by the
# rest of the data rows (left-aligned email column, boolean is_active col).
$ws.Rows("33:33").Insert(-4121, 0)

# Populate the new row with the new user "Ewan Marsh"
$row = 33
$ws.Cells.Item($row, 1).Value = 110032
$ws.Cells.Item($row, 2).Value = 9317596770
$ws.Cells.Item($row, 3).Value = "Ewan Marsh"
$ws.Cells.Item($row, 4).Value = "ewan.marsh@xyz.com"
$ws.Cells.Item($row, 5).Value = 818876433
$ws.Cells.Item($row, 6).Value = "ACT"
$ws.Cells.Item($row, 7).Value = "eng"
$ws.Cells.Item($row, 8).Value = "PWD"
$ws.Cells.Item($row, 9).Value = $true
$ws.Cells.Item($row, 10).Value = "superadmin"
$ws.Cells.Item($row, 11).Value = "now()"

# Update the view: scroll back to top-left and change selection to an
# entire-column selection (L:XFD), matching the saved workbook view.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("L1:XFD1048576").Select()

$wb.Save()
